# Add a "Quantity" column (header + value) to the
# UpdateQuantityAndSaveCart sheet, mirroring the new test data added in
# the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in C1 and its numeric value in C2.
$ws.Range("C1").Value = "Quantity"
$ws.Range("C2").Value = 3

# Reflect the cell the user landed on after entering the data (C2 -> D2).
$ws.Range("D2").Select()
